$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.519.69'
$ws.Range('E2').Value = '  +2.65%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.199.13'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.75'
$ws.Range('E5').Value = '  +2.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '84.46'
$ws.Range('E6').Value = '  +13.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.614'
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  +3.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.01'
$ws.Range('E10').Value = '  +13.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0918'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.13'
$ws.Range('E12').Value = '  +6.01%  '
$ws.Range('E13').Value = '  +2.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.524.45'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.199.12'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.783'
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.473.01'
$ws.Range('E18').Value = '  +2.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000103'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.91'
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.53'
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.38'
$ws.Range('E22').Value = '  +14.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.93'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.88'
$ws.Range('E24').Value = '  -4.05%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.63'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('B27').Value = 'WEMIXToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.41'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.99'
$ws.Range('E28').Value = '  +4.58%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').Value = '  +3.29%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').Value = '  +3.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.34'
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.39'
$ws.Range('E32').Value = '  +2.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0861'
$ws.Range('E33').Value = '  +5.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.35'
$ws.Range('E34').Value = '  +4.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.123'
$ws.Range('E35').Value = '  +2.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.112'
$ws.Range('E36').Value = '  +5.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0359'
$ws.Range('E37').Value = '  +9.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.49'
$ws.Range('E38').Value = '  +7.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.50'
$ws.Range('E39').Value = '  +4.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.84'
$ws.Range('E40').Value = '  +12.19%  '
$ws.Range('E41').Value = '  +2.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '63.26'
$ws.Range('E42').Value = '  +8.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.47'
$ws.Range('E43').Value = '  +6.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.199'
$ws.Range('E44').Value = '  +3.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.07'
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0974'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.24'
$ws.Range('E47').Value = '  +1.33%  '
$ws.Range('E48').Value = '  +5.42%  '
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.439'
$ws.Range('E50').Value = '  -3.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.48'
$ws.Range('E51').Value = '  +15.56%  '
